# --- Rename the weekly-sheet "Requested quantity" header. ---
$wb = $excel.ActiveWorkbook
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- Rename the monthly-sheet "Requested quantity" header. ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end of the workbook. ---
# Clone the "Weekly Quantity" sheet so the new sheet inherits identical
# sheet properties / page setup / cell styles (bold+bordered header,
# date-time number format on column A), then wipe its old data.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsWeekly.Copy($null, $lastSheet)
$wsForecast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast.Name = "PO Forecast"
$wsForecast.Range("A1:B52").ClearContents()

# Extend the header's bold/bordered style to columns C and D.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("C1:D1").PasteSpecial(-4122)

# Extend the date-time number format down column A for the extra rows
# beyond what "Weekly Quantity" originally had (it only had 52 rows).
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A53:A60").PasteSpecial(-4122)

# Header row.
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper).
$forecastData = @(
    @(45053.99999999999, 10, -29.78959112350789, 54.68011390997718),
    @(45060.99999999999, 11, -30.60969324237815, 56.08263424599006),
    @(45067.99999999999, 11, -35.88291899272289, 53.4295530251608),
    @(45074.99999999999, 12, -33.1661735158172, 53.0595262742677),
    @(45081.99999999999, 12, -31.79284613188408, 53.72193107229371),
    @(45088.99999999999, 12, -29.81833193112677, 55.97134536744654),
    @(45095.99999999999, 13, -30.13309829688218, 56.83042209524797),
    @(45102.99999999999, 13, -30.74331872551251, 55.62075073756303),
    @(45109.99999999999, 14, -34.12259104570421, 56.78847434103568),
    @(45116.99999999999, 14, -26.60110110399815, 57.59152668021489),
    @(45123.99999999999, 15, -31.16937453715015, 57.39145674041921),
    @(45137.99999999999, 16, -28.74734632368287, 58.30528319784066),
    @(45144.99999999999, 16, -30.22728915196477, 59.76109997065775),
    @(45151.99999999999, 16, -29.70489128261889, 61.27191607189053),
    @(45158.99999999999, 17, -26.26206256426359, 60.53941025553286),
    @(45165.99999999999, 17, -27.37205913657696, 59.97675185555628),
    @(45172.99999999999, 18, -27.80091171381772, 62.99575160415458),
    @(45179.99999999999, 18, -26.27200556934087, 62.49482968637154),
    @(45186.99999999999, 19, -24.45590138803792, 61.08498550352128),
    @(45193.99999999999, 19, -24.6351863071521, 62.78002480688625),
    @(45200.99999999999, 20, -24.11595424077906, 65.91559068211774),
    @(45207.99999999999, 20, -22.1442469582629, 63.3334910206969),
    @(45214.99999999999, 21, -22.38547192573702, 63.8954810225739),
    @(45228.99999999999, 21, -23.39779095752704, 66.21880296473425),
    @(45235.99999999999, 22, -21.15312838064596, 69.65608495665336),
    @(45242.99999999999, 22, -22.01593175148977, 66.22360940223946),
    @(45249.99999999999, 23, -23.02234950688787, 66.09684548879912),
    @(45256.99999999999, 23, -21.21448590051637, 67.63453617766719),
    @(45312.99999999999, 27, -15.45712907792346, 72.11421986483363),
    @(45319.99999999999, 27, -17.1107397297866, 70.46217261784763),
    @(45326.99999999999, 28, -15.15895961083634, 72.78698017474913),
    @(45333.99999999999, 28, -15.50323136520536, 72.33219483437539),
    @(45340.99999999999, 29, -13.27969291337733, 71.9125305671738),
    @(45347.99999999999, 29, -15.14668807640055, 72.75196308841775),
    @(45354.99999999999, 29, -14.08107954422655, 69.28190568126253),
    @(45361.99999999999, 30, -16.79624109244013, 73.95339070302418),
    @(45368.99999999999, 30, -16.17583525117433, 71.61970908457141),
    @(45375.99999999999, 31, -15.48640805852773, 74.0295249105888),
    @(45382.99999999999, 31, -12.52010790807461, 74.61605693310997),
    @(45389.99999999999, 32, -10.85157376333509, 76.74384975342663),
    @(45396.99999999999, 32, -11.04331396435813, 79.01247210351688),
    @(45403.99999999999, 33, -10.94493346998989, 76.84733037723193),
    @(45410.99999999999, 33, -10.17647632179772, 72.54017071282136),
    @(45417.99999999999, 34, -12.62541195948917, 74.98421134912077),
    @(45424.99999999999, 34, -8.807793555575181, 77.25446515789871),
    @(45431.99999999999, 34, -10.33809694247973, 75.89769130327286),
    @(45438.99999999999, 35, -8.644224168870087, 79.45974535177487),
    @(45445.99999999999, 35, -9.010205491337505, 81.65401817841629),
    @(45452.99999999999, 36, -6.410824218740843, 79.1698354426543),
    @(45459.99999999999, 36, -7.686875639074508, 80.59992458408334),
    @(45466.99999999999, 37, -9.906240247935699, 80.07371179838127),
    @(45473.99999999999, 37, -6.214326949200696, 78.21744166387434),
    @(45480.99999999999, 38, -6.862661152833891, 81.48020425172669),
    @(45487.99999999999, 38, -7.374173965848462, 81.69446317855009),
    @(45494.99999999999, 38, -4.396600714669391, 84.7534081470879),
    @(45501.99999999999, 39, -6.856681941696261, 78.14873647541305),
    @(45508.99999999999, 39, -4.371743170368507, 86.30436774986202),
    @(45515.99999999999, 40, -4.635537438395793, 83.08473343006045),
    @(45522.99999999999, 40, -3.132280279917185, 82.86909776567772)
)

$r = 2
foreach ($row in $forecastData) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

Write-Output "PO Forecast sheet added"
